# TC_62.xlsx edit script
# - Rename sheet "My Series" -> "Data"
# - Update the hidden CEIC add-in metadata blob stored in cell A1's comment
# - Change the numeric format used by B13:B16 from "0.000" to "###0.000"
# - Change A11's label from "Function Description" to "Function Information"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the worksheet tab
$ws.Name = "Data"

# 2. Update the comment attached to A1 with the new compressed metadata payload
$newCommentText = "QhwAAB+LCAAAAAAAAAOlGdtu48b1Vwg9tUAlkpLt9RpjBrrZESpbhiTX67wUI3JssSaHKmdoW28pkCBFmqIoik2RXtGnFAW6XbQJkO728i/B2rt96i/0zIU3Sc5a7sJYc85tzpw51zF65zoMjEsSMz+iuxW7ZlUMQt3I8+n5biXhZ1V7q/KOg7rXLgmOcIxDwoHYAC7Kdq6Zv1uZcj7bMc2rq6vaVaMWxedm3bJs88lBf+ROSYirPmUcU5dUMi7v7VwVB7W98IBw7GGOFedupTfq1drEdzsAO8AUn5O41kqYTwljXcp97hMmOGOCOWl3Dr6nDubUa1s1G5lL8JyylfiBp+hKlAqu6WBbMvZD4tQte7tqbVcb9ti2dja2dxpbtUdW/b2UMSNEfcz4iMSXvisBI47DmWS3thu2bW1sN7aQuZIIZOUGcNAg8Ibk0mfEa5MgYGtZxNQX2HQ5nHo9Y1rILPBqQQ9XYT/Gs+nY5wFZT43hQcsIqdYlF+KgvSgmLtjvQSodkqtBrM06nvUBO576MZ938HxtWceMxIOZMNJ6rA7qRJQ3AxLz4xncNfHAFQDh8DghyLwDmTN1fObCt08T4jlnOGBFphISnUTxBZthlxxCHJtCxhUNIuyBw3Gfcd9luYAlDDqKoxmIhN1bUeDtgVit4gpEJrpHwcZi31YUXSxKLyORvFZ5wXCpIeYp+RIcjabR1YAG81EyYW7sT4jXaaXUK3FIRKTmbieMRyFokYOQghUgc/gHEbgIRh3i+iEOjgKwInMaIKUEQM2ER2c+b0dBEtLMnAtQdAInGpPr1IDZEg3gcqmweUR7NCVXRCtRZYZhdJVtuYyQNiiAm8xNfWwZsUjcAVh6e8sYeSHikHt+APWheBUFaNkpRlNC+EqPUBgkUuGeqDhOay72RGYOQeCY4N0AhWRqWVXLhp+xZe3IH9g5Q6Mu9VK6jSW6FIkOk3AwgQi+lGdybMAtgBCcImgFmF4A9MTn08Nmqv0KDFJnvpN+GYcgVmcBnktwZpciDPWoGyQeUSmgR8+kTwrd1DXeiUZLoD4EtYMwnY/nM8jEzN/h8LFbgdq8w3gM1b/iuFFCeTwXuQKZmvRtPCyZULkBDu7NcxaTHybQdMz3Euq2I+/+u3nKOsfU5/fXMEpilQDvzyKtJ1JhwjpEJBWZ5u/N765zJhavRR5SEkbUd+9vbTCy0N57wEFYGlUFDvMO3+qQM5wE0NZwqD7neWFZAKMmu1ikKYLQcRykycERTSODrtH1wpoLdVV0RjU3CgXAhGbtZITMIr1oDlzSped9TM8TKL9ZAC7Cs9QkKsc4xpSJ42TVdiFLrSZCaUCrLsBRUT5IpMVUlEeAReYCHRqTcBbFODgAw/h7+n50KwFV+gDzqV5B2g+ImxrZzFkzrrJmqeJvI5P5Wx1DRIbOJwtASSTOovrTnCaHIXHKA/DfoI0DfxKr9KOVWImDC8v7pjRRicOt2UOldwAzChSm75K5aFzzhYZLl7VThHJgkXGc0XBju75pNepQ88UayRMPCQ6MLng9J0aPXhLGQ2DbMYaE+R58+TjYMd4lE+JDtZAm0vl6be4iH9pLE6JUpQn6liFlAijF5z7k22XCDJMzOKcEx8G8QKiO2o9coLv98b9vfvP81YvPbz9++ubLD/77j1+9+ufPbp59CB+3f/3bzSe/VMdUxGiMJwGRCo1b29tWYwP8LAMhYVxTNo1e4nIJOz2VvWK2RnrmkYt2t9fe77dkPsmAKbvKvaYYp+ZRki9H6hByI3mlZuoJisQZp/lJr0vYQi53xHhzScrURfxdjMoWr19+/vrln+/k1gbLmxL78ePNql2/R89i20t0Wc8ixkY1FaRDKTQ4m9V6vUC8QIOGMFvDEJDZqec5Ddt6bNUblp3lci9z5FVEiygtaYzPzQU+BWqrPiJzgeI6RUrHH0OIZGgVCoWFdtEvf/LmL09LVNq6GlKWAsrJei82M9OFFH04HBujwfGw3TXG3ZHwkxxXoFPCv4FY757FU8mpKE1w8B2DUPGeYlRgVqgY0ZlBsDs15hCJhTgsOdsqqNrogSIXtdyPo2SmbqTAkENXUGbZZCXHilwjcdKeS0knR60gV7re/P2LVQz6IJ2888ueBoowVMIoUAGvo/azf7366qNXL17cPv/5zVc/KknQ+2QTM/g5RFNxmbk9pDxdbxYg6GQkjXlhfb9QXzRQTBxHkU85c+wtOWzoFQJWW0iTv1EvhJInBUt7AXwBgt7FrHvNdWA7h8gsA0DPGYZqG+VjWQZQOTy3639++7vbX39x++nzNx/96ebjP9588unrl79/8+wPKupunz6//ekzneUXC4HURQx7qgk05MuBa4hoNETtNr5+/xcGjbgBLYeRyIz09fufFYQJRWVzkkuGli5TpKzCEmmRWfAZBVUyHUp8GYtqANqihDUyCl3Eopnv5pu8VxWiRNxJxLd642rCiBFBN/VtOEmZOGe+L59mUSX16JFVt+saq7QRR5hgVjD9fhBNoMlIEXI2XyApcX0zQ04r99vvD1rNfk6ilBjEHozylnhgEB8obSlFSemxdJW6WgECWGj83CQQrylLZMuoTHIhjZn6XeKs6Yn058Doby/N9CUK1E7iWDVEVD9zj5IZNMPp28vdePmiV+h/D1WvWuyI83WvU8bDuoCFQlhGC4DEy9SkUSpN9Zh4AlHt7KEwTb4EXOkVEMyhX7JVp3UJfWVsirzTjeMoXpl8ckxKdgCdNGQUM7d4RiPvVHXdXn5XKSBNeNmHmvz0CaMOCQhf75nXzLkPossH88Ldr8vaY4PA08Zcb/TIzJILKL51C0f5f5+6lbM14xgaK/E2tvbbdDq4DmHeXVMbdRTJKCZA2F0/Hu/5MeNPRCbQXwpymkFOVYf6RAxc6kOuTx17UwGAwCxKN0tqpqHL1V8coqDvh/6aY6GVxndZCNhyNlMtXG89TxGl5ZBcQ4NZkABJcfIDKBvq5WQdacphIZdm/OKRj/nnU76uYo8mmHhkYlXdCalXNzxru/qYkEbVtuF/7NbrlrUpngi1cMgcPrlacxMzvbD8L4HO/wA0MGnoQhwAAA=="

$comment = $ws.Comments.Item(1)
$comment.Text($newCommentText)

# 3. Change the number format for the data column (was "0.000", now "###0.000")
$ws.Range("B13:B16").NumberFormat = "###0.000"

# 4. Relabel the "Function Description" row
$ws.Range("A11").Value = "Function Information"
